# Generate Report for Handback
# The source file "0c05a52b-0c83-4586-9214-fbdd36cbc2b1.md" has been handed
# back and is in sync with en-US, so update the status/handback info for
# that row across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_dede = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: row 2 is the 0c05a52b-... file ---
$ws_overview.Range("B2").Value = $statusText
$ws_overview.Range("C2").Value = $statusText

# --- zh-cn sheet: row 2 is the 0c05a52b-... file ---
$ws_zhcn.Range("C2").Value = $statusText
$ws_zhcn.Range("F2").Value = "0c05a52b-0c83-4586-9214-fbdd36cbc2b1.md"
$ws_zhcn.Range("G2").Value = "0c05a52b-0c83-4586-9214-fbdd36cbc2b1.5ecd15e3fdae5f448c36f40c571c9b41a8f26c74.zh-cn.xlf"
$ws_zhcn.Range("H2").Value = "2016-03-19 08:30:39"

$ws_zhcn.Hyperlinks.Add(
    $ws_zhcn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f3f0bcc76a05ce10f5c256e5d9009383cc543374/e2e/0c05a52b-0c83-4586-9214-fbdd36cbc2b1.md",
    "",
    "",
    "0c05a52b-0c83-4586-9214-fbdd36cbc2b1.md"
) | Out-Null
$ws_zhcn.Hyperlinks.Add(
    $ws_zhcn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/086b6d37dc6c43a7cf9719a91987f1bcaf05af36/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0c05a52b-0c83-4586-9214-fbdd36cbc2b1.5ecd15e3fdae5f448c36f40c571c9b41a8f26c74.zh-cn.xlf",
    "",
    "",
    "0c05a52b-0c83-4586-9214-fbdd36cbc2b1.5ecd15e3fdae5f448c36f40c571c9b41a8f26c74.zh-cn.xlf"
) | Out-Null

# --- de-de sheet: row 2 is the 0c05a52b-... file ---
$ws_dede.Range("C2").Value = $statusText
$ws_dede.Range("F2").Value = "0c05a52b-0c83-4586-9214-fbdd36cbc2b1.md"
$ws_dede.Range("G2").Value = "0c05a52b-0c83-4586-9214-fbdd36cbc2b1.5ecd15e3fdae5f448c36f40c571c9b41a8f26c74.de-de.xlf"
$ws_dede.Range("H2").Value = "2016-03-19 08:30:44"

$ws_dede.Hyperlinks.Add(
    $ws_dede.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f3f0bcc76a05ce10f5c256e5d9009383cc543374/e2e/0c05a52b-0c83-4586-9214-fbdd36cbc2b1.md",
    "",
    "",
    "0c05a52b-0c83-4586-9214-fbdd36cbc2b1.md"
) | Out-Null
$ws_dede.Hyperlinks.Add(
    $ws_dede.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a10bab7a56a4a10287c658d1e230a7c0ac061b95/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0c05a52b-0c83-4586-9214-fbdd36cbc2b1.5ecd15e3fdae5f448c36f40c571c9b41a8f26c74.de-de.xlf",
    "",
    "",
    "0c05a52b-0c83-4586-9214-fbdd36cbc2b1.5ecd15e3fdae5f448c36f40c571c9b41a8f26c74.de-de.xlf"
) | Out-Null

Write-Host "Handback report generated."
